# :sparkles: #2 컬럼 auto width 설정 옵션 추가
#
# Adds three date/time columns (date, date2, dateTime) to every sheet,
# refreshes the existing numeric sample data, freezes the header row and
# sets explicit (auto-fit-derived) column widths on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Per-sheet data: existing C/D/E values to overwrite + new F/G/H values
# and the (COM ColumnWidth-space) widths for columns A-H.
# ---------------------------------------------------------------------
$sheetData = @{
    1 = @{
        Rows = @(
            @{ Row = 2; C = 41385.0;  D = 295738.875;   E = 208787.0; F = 45020.0; G = 45020.69528284722; H = 45020.69528284722 },
            @{ Row = 3; C = 306765.0; D = 455219.03125;  E = 715353.0; F = 45020.0; G = 45020.6952828588;  H = 45020.6952828588  }
        )
        Widths = @(2.5, 14.666666666666666, 7.5, 10.166666666666666, 8.5, 10.5, 18.833333333333332, 18.833333333333332)
    }
    2 = @{
        Rows = @(
            @{ Row = 2; C = 546439.0; D = 345885.5625;  E = 480987.0; F = 45020.0; G = 45020.6952828588; H = 45020.6952828588 },
            @{ Row = 3; C = 150976.0; D = 113105.0625;  E = 490072.0; F = 45020.0; G = 45020.6952828588; H = 45020.6952828588 }
        )
        Widths = @(2.5, 22.833333333333332, 7.5, 10.166666666666666, 8.5, 10.5, 18.833333333333332, 18.833333333333332)
    }
    3 = @{
        Rows = @(
            @{ Row = 2; C = 928855.0; D = 179729.28125; E = 242430.0; F = 45020.0; G = 45020.6952828588; H = 45020.6952828588 }
        )
        Widths = @(2.5, 21.833333333333332, 7.5, 10.166666666666666, 8.5, 10.5, 18.833333333333332, 18.833333333333332)
    }
}

for ($sheetIdx = 1; $sheetIdx -le 3; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $info = $sheetData[$sheetIdx]
    $lastRow = ($info.Rows | Measure-Object -Property Row -Maximum).Maximum

    # --- Header row: new F1/G1/H1 cells, formatted like the existing E1 header ---
    $ws.Range("E1").Copy()
    $ws.Range("F1:H1").PasteSpecial(-4122)
    $ws.Range("F1").Value = "date"
    $ws.Range("G1").Value = "date2"
    $ws.Range("H1").Value = "dateTime"

    # --- New data cells F:H, formatted like the existing E column data cells ---
    $ws.Range("E2:E" + $lastRow).Copy()
    $ws.Range("F2:H" + $lastRow).PasteSpecial(-4122)

    foreach ($r in $info.Rows) {
        $ws.Range("C" + $r.Row).Value = $r.C
        $ws.Range("D" + $r.Row).Value = $r.D
        $ws.Range("E" + $r.Row).Value = $r.E
        $ws.Range("F" + $r.Row).Value = $r.F
        $ws.Range("G" + $r.Row).Value = $r.G
        $ws.Range("H" + $r.Row).Value = $r.H
    }

    # --- Number formats: F = date only, G/H = date + time ---
    $ws.Range("F2:F" + $lastRow).NumberFormat = "yyyy-mm-dd"
    $ws.Range("G2:H" + $lastRow).NumberFormat = "yyyy-mm-dd hh:mm:ss"

    # --- Auto-width columns A-H ---
    $widths = $info.Widths
    for ($c = 1; $c -le 8; $c++) {
        $ws.Columns.Item($c).ColumnWidth = $widths[$c - 1]
    }

    # --- Freeze the header row ---
    $ws.Activate()
    $ws.Range("A2").Select()
    $excel.ActiveWindow.FreezePanes = $true
}

$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()
